$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 295) holds a "last changed" date serial value
# that was bumped by one day (2023-09-12 -> 2023-09-13, serial 45181 -> 45182)
# for every data row in the sheet.
$ws.Range("C2:C295").Value = 45182
